$wb = $excel.ActiveWorkbook

# --- 1. Insert a new "Player Info" sheet before the existing "ODI Batting" sheet ---
$wsBattingTemp = $wb.Worksheets.Item("ODI Batting")
$wsInfo = $wb.Worksheets.Add($wsBattingTemp)
$wsInfo.Name = "Player Info"

# Re-fetch stable (name-based) references now that the sheet order has shifted
$wsBatting = $wb.Worksheets.Item("ODI Batting")
$wsBowling = $wb.Worksheets.Item("ODI Bowling")

# Header row formatting (bold, thin border all around, centered/top aligned)
$wsInfo.Range("A1:D1").Font.Bold = $true
$wsInfo.Range("A1:D1").HorizontalAlignment = -4108
$wsInfo.Range("A1:D1").VerticalAlignment = -4160
$wsInfo.Range("A1:D1").Borders.LineStyle = 1

$wsInfo.Range("A1").Value = "ID"
$wsInfo.Range("B1").Value = "NAME"
$wsInfo.Range("C1").Value = "BATTING_HAND"
$wsInfo.Range("D1").Value = "BOWL_STYLE"

$wsInfo.Range("A2").NumberFormat = "@"
$wsInfo.Range("A2").Value = "4684"
$wsInfo.Range("B2").Value = "Washington Sundar"
$wsInfo.Range("C2").Value = "Left Handed"
$wsInfo.Range("D2").Value = "Right Arm Off Break"

# --- 2. "ODI Batting": rename MATCH_CARD_LINK (col D) -> MATCH_CODE and
#        replace the full scorecard URL with just the trailing numeric match code ---
$wsBatting.Cells.Item(1, 4).Value = "MATCH_CODE"
$lastRowBatting = $wsBatting.UsedRange.Rows.Count
for ($i = 2; $i -le $lastRowBatting; $i++) {
    $cell = $wsBatting.Cells.Item($i, 4)
    $txt = $cell.Text
    if ($txt -ne "") {
        $parts = $txt.Split("=")
        $code = $parts[$parts.Length - 1]
        $cell.NumberFormat = "@"
        $cell.Value = $code
    }
}

# --- 3. "ODI Bowling": rename MATCH_CARD_LINK (col B) -> MATCH_CODE and
#        replace the full scorecard URL with just the trailing numeric match code ---
$wsBowling.Cells.Item(1, 2).Value = "MATCH_CODE"
$lastRowBowling = $wsBowling.UsedRange.Rows.Count
for ($i = 2; $i -le $lastRowBowling; $i++) {
    $cell = $wsBowling.Cells.Item($i, 2)
    $txt = $cell.Text
    if ($txt -ne "") {
        $parts = $txt.Split("=")
        $code = $parts[$parts.Length - 1]
        $cell.NumberFormat = "@"
        $cell.Value = $code
    }
}
